# Generate Report for Handback
# This script regenerates the localization-status report after a new
# handback took place: the status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps are
# refreshed, and the stale "handback not latest" error is cleared now
# that everything is back in sync.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column width (character units) that corresponds to the OOXML stored
# width of ~29.9777047293527 (status column growing to fit the new,
# longer status text).
$wideColWidth = 29.144371396019366
# Column width (character units) that corresponds to the OOXML stored
# width of ~13.7470528738839 (error column shrinking back down now that
# the long error text has been cleared).
$narrowColWidth = 12.913719540550566

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $wideColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-13 17:07:13"
$wsZhCn.Range("P2").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = $wideColWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $narrowColWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-13 17:07:24"
$wsDeDe.Range("P2").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = $wideColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $narrowColWidth
